$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 498.6
$ws.Range("J33").Value = 175
$ws.Range("L33").Value = 175
$ws.Range("N33").Value = -633
$ws.Range("H40").Value = 2838.2104
$ws.Range("I40").Value = 1880.4166
$ws.Range("J40").Value = 4480.143
$ws.Range("K40").Value = 1880.4166
$ws.Range("L40").Value = 4480.143
$ws.Range("M40").Value = -1705.4166
$ws.Range("N40").Value = -4830.143
$ws.Range("H74").Value = 8803.5
$ws.Range("I74").Value = 6178.4287
$ws.Range("J74").Value = 11428.571
$ws.Range("K74").Value = 6178.4287
$ws.Range("L74").Value = 11428.571
$ws.Range("M74").Value = -5242.4287
$ws.Range("N74").Value = -13300.571
$ws.Range("H77").Value = 8803.5
$ws.Range("I77").Value = 6178.4287
$ws.Range("J77").Value = 11428.571
$ws.Range("K77").Value = 30892.1435
$ws.Range("L77").Value = 57142.855
$ws.Range("M77").Value = -26212.1435
$ws.Range("N77").Value = -66502.855
$ws.Range("H80").Value = 2173.5134
$ws.Range("I80").Value = 1495.125
$ws.Range("J80").Value = 2690.3809
$ws.Range("K80").Value = 4485.375
$ws.Range("L80").Value = 8071.1427
$ws.Range("M80").Value = -3487.375
$ws.Range("N80").Value = -10067.1427
$ws.Range("H83").Value = 2173.5134
$ws.Range("I83").Value = 1495.125
$ws.Range("J83").Value = 2690.3809
$ws.Range("K83").Value = 13456.125
$ws.Range("L83").Value = 24213.4281
$ws.Range("M83").Value = -8464.125
$ws.Range("N83").Value = -34197.4281
$ws.Range("H96").Value = 244.5
$ws.Range("I96").Value = 244.5
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 733.5
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 639.5
$ws.Range("N96").ClearContents()
$ws.Range("H103").Value = 967.8
$ws.Range("I103").Value = 946.6667
$ws.Range("J103").Value = 999.5
$ws.Range("K103").Value = 2840.0001
$ws.Range("L103").Value = 2998.5
$ws.Range("M103").Value = -2254.0001
$ws.Range("N103").Value = -4170.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16686.85
$ws.Range("I2").Value = 17301.947
$ws.Range("K2").Value = 17301.947
$ws.Range("M2").Value = -17188.947
$ws.Range("H32").Value = 8909.666999999999
$ws.Range("I32").Value = 7039.3584
$ws.Range("J32").Value = 16534.77
$ws.Range("K32").Value = 7039.3584
$ws.Range("L32").Value = 16534.77
$ws.Range("M32").Value = -6752.3584
$ws.Range("N32").Value = -17108.77
$ws.Range("H51").Value = 38495
$ws.Range("J51").Value = 38495
$ws.Range("L51").Value = 38495
$ws.Range("N51").Value = -40007
$ws.Range("H116").Value = 16686.85
$ws.Range("I116").Value = 17301.947
$ws.Range("K116").Value = 17301.947
$ws.Range("M116").Value = -15007.947
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16686.85
$ws.Range("I3").Value = 17301.947
$ws.Range("K3").Value = 17301.947
$ws.Range("M3").Value = -17187.947
$ws.Range("H35").Value = 200000
$ws.Range("I35").Value = 200000
$ws.Range("K35").Value = 200000
$ws.Range("M35").Value = -199690
$ws.Range("H86").Value = 2652.577
$ws.Range("J86").Value = 5061.125
$ws.Range("L86").Value = 5061.125
$ws.Range("N86").Value = -7307.125
$ws.Range("H89").Value = 2652.577
$ws.Range("J89").Value = 5061.125
$ws.Range("L89").Value = 25305.625
$ws.Range("N89").Value = -36537.625
$ws.Range("H105").Value = 3570
$ws.Range("I105").Value = 4093.3333
$ws.Range("K105").Value = 4093.3333
$ws.Range("M105").Value = -2346.3333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4179.1177
$ws.Range("I31").Value = 1192
$ws.Range("K31").Value = 1192
$ws.Range("M31").Value = -897
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H34").Value = 4179.1177
$ws.Range("I34").Value = 1192
$ws.Range("K34").Value = 1192
$ws.Range("M34").Value = -990
$ws.Range("H58").Value = 1562.875
$ws.Range("J58").Value = 3171.125
$ws.Range("L58").Value = 3171.125
$ws.Range("N58").Value = -3577.125
$ws.Range("H86").Value = 41472.832
$ws.Range("I86").Value = 55709
$ws.Range("J86").Value = 21542.2
$ws.Range("K86").Value = 55709
$ws.Range("L86").Value = 21542.2
$ws.Range("M86").Value = -54586
$ws.Range("N86").Value = -23788.2
$ws.Range("H89").Value = 41472.832
$ws.Range("I89").Value = 55709
$ws.Range("J89").Value = 21542.2
$ws.Range("K89").Value = 278545
$ws.Range("L89").Value = 107711
$ws.Range("M89").Value = -272929
$ws.Range("N89").Value = -118943
$ws.Range("H99").Value = 13206139
$ws.Range("J99").Value = 25008876
$ws.Range("L99").Value = 25008876
$ws.Range("N99").Value = -25011872
$ws.Range("H105").Value = 13894589
$ws.Range("I105").Value = 1640.4
$ws.Range("J105").Value = 83359336
$ws.Range("K105").Value = 1640.4
$ws.Range("L105").Value = 83359336
$ws.Range("M105").Value = 106.5999999999999
$ws.Range("N105").Value = -83362830
$ws.Range("H126").Value = 13206139
$ws.Range("J126").Value = 25008876
$ws.Range("L126").Value = 75026628
$ws.Range("N126").Value = -75031568
$ws.Range("H136").Value = 1562.875
$ws.Range("J136").Value = 3171.125
$ws.Range("L136").Value = 9513.375
$ws.Range("N136").Value = -14613.375

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 36.61111
$ws.Range("I2").Value = 8.666667
$ws.Range("J2").Value = 64.55556
$ws.Range("K2").Value = 52.000002
$ws.Range("L2").Value = 387.33336
$ws.Range("M2").Value = 60.999998
$ws.Range("N2").Value = -613.33336
$ws.Range("H14").Value = 248
$ws.Range("I14").Value = 248
$ws.Range("K14").Value = 744
$ws.Range("M14").Value = -571
$ws.Range("H80").Value = 4239.5
$ws.Range("J80").Value = 4538.3335
$ws.Range("L80").Value = 13615.0005
$ws.Range("N80").Value = -15487.0005
$ws.Range("H83").Value = 4239.5
$ws.Range("J83").Value = 4538.3335
$ws.Range("L83").Value = 40845.0015
$ws.Range("N83").Value = -50205.0015
$ws.Range("H140").Value = 1397.7646
$ws.Range("I140").Value = 1397.7646
$ws.Range("K140").Value = 4193.293799999999
$ws.Range("M140").Value = 986.7062000000005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 523.7692
$ws.Range("I97").Value = 534.4545000000001
$ws.Range("K97").Value = 534.4545000000001
$ws.Range("M97").Value = -38.45450000000005
$ws.Range("H113").Value = 5737.273
$ws.Range("I113").Value = 3551.8333
$ws.Range("J113").Value = 8359.799999999999
$ws.Range("K113").Value = 3551.8333
$ws.Range("L113").Value = 8359.799999999999
$ws.Range("M113").Value = -1381.8333
$ws.Range("N113").Value = -12699.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 10000
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("N38").Value = -10820
$ws.Range("H74").Value = 34966.668
$ws.Range("J74").Value = 34966.668
$ws.Range("L74").Value = 34966.668
$ws.Range("N74").Value = -36962.668
$ws.Range("H77").Value = 34966.668
$ws.Range("J77").Value = 34966.668
$ws.Range("L77").Value = 104900.004
$ws.Range("N77").Value = -114884.004
$ws.Range("H100").Value = 74003
$ws.Range("I100").Value = 188333
$ws.Range("J100").Value = 5405
$ws.Range("K100").Value = 188333
$ws.Range("L100").Value = 5405
$ws.Range("M100").Value = -187792
$ws.Range("N100").Value = -6487
$ws.Range("H132").Value = 5472.1333
$ws.Range("I132").Value = 3025.6
$ws.Range("K132").Value = 9076.799999999999
$ws.Range("M132").Value = -6546.799999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6099.1816
$ws.Range("I81").Value = 7945.75
$ws.Range("K81").Value = 15891.5
$ws.Range("M81").Value = -14830.5
$ws.Range("H84").Value = 6099.1816
$ws.Range("I84").Value = 7945.75
$ws.Range("K84").Value = 79457.5
$ws.Range("M84").Value = -74153.5
$ws.Range("H96").Value = 2197.75
$ws.Range("I96").Value = 1945.5
$ws.Range("J96").Value = 2450
$ws.Range("K96").Value = 1945.5
$ws.Range("L96").Value = 2450
$ws.Range("M96").Value = -572.5
$ws.Range("N96").Value = -5196
$ws.Range("H100").Value = 1133.5
$ws.Range("I100").Value = 1160.2
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 2320.4
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1779.4
$ws.Range("N100").Value = -3082
$ws.Range("H113").Value = 387.85715
$ws.Range("I113").Value = 283.33334
$ws.Range("J113").Value = 466.25
$ws.Range("K113").Value = 850.0000200000001
$ws.Range("L113").Value = 1398.75
$ws.Range("M113").Value = 1319.99998
$ws.Range("N113").Value = -5738.75
$ws.Range("H132").Value = 1971.6
$ws.Range("I132").Value = 1610
$ws.Range("J132").Value = 3116.6667
$ws.Range("K132").Value = 4830
$ws.Range("L132").Value = 9350.000100000001
$ws.Range("M132").Value = -2300
$ws.Range("N132").Value = -14410.0001
